$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 3 de Septiembre de 2020 a las 08:19"

# --- Plain data refresh for countries whose rank/position did not move ---

# Row 6 (India)
$ws.Cells.Item(6,2).Value = 3853406
$ws.Cells.Item(6,3).Value = 4438
$ws.Cells.Item(6,4).Value = 2970492
$ws.Cells.Item(6,5).Value = 815428

# Row 19 (Pakistan)
$ws.Cells.Item(19,2).Value = 297014
$ws.Cells.Item(19,3).Value = 424
$ws.Cells.Item(19,4).Value = 281925
$ws.Cells.Item(19,5).Value = 8761
$ws.Cells.Item(19,7).Value = 10
$ws.Cells.Item(19,8).Value = 6328

# Row 29 (Israel)
$ws.Cells.Item(29,2).Value = 122539
$ws.Cells.Item(29,3).Value = 1075
$ws.Cells.Item(29,4).Value = 97872
$ws.Cells.Item(29,5).Value = 23698

# Row 33 (Kazajistan)
$ws.Cells.Item(33,4).Value = 97967
$ws.Cells.Item(33,5).Value = 6477

# --- Armenia / Kirguistan swap places (rows 59-60) + Kirguistan's data refresh ---
# Row 59 becomes Kirguistan (new data); row 60 becomes Armenia (its old, unchanged data)
$ws.Cells.Item(59,1).Value = "Kirguistan"
$ws.Cells.Item(59,2).Value = 44135
$ws.Cells.Item(59,3).Value = 99
$ws.Cells.Item(59,4).Value = 39174
$ws.Cells.Item(59,5).Value = 3902
$ws.Cells.Item(59,8).Value = 1059

# Row 60
$ws.Cells.Item(60,1).Value = "Armenia"
$ws.Cells.Item(60,2).Value = 44075
$ws.Cells.Item(60,4).Value = 38631
$ws.Cells.Item(60,5).Value = 4560
$ws.Cells.Item(60,8).Value = 884

# Row 72 (Australia)
$ws.Cells.Item(72,4).Value = 21912
$ws.Cells.Item(72,5).Value = 3459

# Row 73 (El Salvador)
$ws.Cells.Item(73,4).Value = 14904
$ws.Cells.Item(73,5).Value = 10261
$ws.Cells.Item(73,7).Value = 8
$ws.Cells.Item(73,8).Value = 739

# Row 124 (Tailandia)
$ws.Cells.Item(124,2).Value = 3427
$ws.Cells.Item(124,3).Value = 2
$ws.Cells.Item(124,4).Value = 3277
$ws.Cells.Item(124,5).Value = 92

# --- Jamaica jumps ahead of Mali/Angola (rows 133-135) with a data refresh; ---
# --- Mali and Angola shift down one row each, keeping their own old data   ---
# Row 133 becomes Jamaica (new data)
$ws.Cells.Item(133,1).Value = "Jamaica"
$ws.Cells.Item(133,2).Value = 2822
$ws.Cells.Item(133,3).Value = 139
$ws.Cells.Item(133,4).Value = 900
$ws.Cells.Item(133,5).Value = 1895
$ws.Cells.Item(133,7).Value = 3
$ws.Cells.Item(133,8).Value = 27

# Row 134 becomes Mali (= old row 133 data, unchanged values, just shifted down)
$ws.Cells.Item(134,1).Value = "Mali"
$ws.Cells.Item(134,2).Value = 2802
$ws.Cells.Item(134,4).Value = 2185
$ws.Cells.Item(134,5).Value = 491
$ws.Cells.Item(134,8).Value = 126

# Row 135 becomes Angola (= old row 134 data, unchanged values, just shifted down)
$ws.Cells.Item(135,1).Value = "Angola"
$ws.Cells.Item(135,2).Value = 2777
$ws.Cells.Item(135,4).Value = 1115
$ws.Cells.Item(135,5).Value = 1550
$ws.Cells.Item(135,8).Value = 112

# --- Birmania jumps ahead of Vietnam/Chad (rows 164-166) with a data refresh; ---
# --- Vietnam and Chad shift down one row each, keeping their own old data     ---
# Row 164 becomes Birmania (new data)
$ws.Cells.Item(164,1).Value = "Birmania"
$ws.Cells.Item(164,2).Value = 1058
$ws.Cells.Item(164,3).Value = 63
$ws.Cells.Item(164,4).Value = 359
$ws.Cells.Item(164,5).Value = 693
$ws.Cells.Item(164,8).Value = 6

# Row 165 becomes Vietnam (= old row 164 data, unchanged values, just shifted down)
$ws.Cells.Item(165,1).Value = "Vietnam"
$ws.Cells.Item(165,2).Value = 1046
$ws.Cells.Item(165,4).Value = 746
$ws.Cells.Item(165,5).Value = 266
$ws.Cells.Item(165,8).Value = 34

# Row 166 becomes Republica del Chad (= old row 165 data, unchanged values, just shifted down)
$ws.Cells.Item(166,1).Value = "Republica del Chad"
$ws.Cells.Item(166,2).Value = 1017
$ws.Cells.Item(166,4).Value = 904
$ws.Cells.Item(166,5).Value = 36
$ws.Cells.Item(166,8).Value = 77

# Row 187 (Butan)
$ws.Cells.Item(187,4).Value = 144
$ws.Cells.Item(187,5).Value = 83
